$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{A="ECs"; D="ECs"; E=3; F=1; G=29.947775; H=89.84332500000001; I=0.9303126840830549; J=0.930312684083055; K=2; L=0.6666666666666666; M=10.53177933333333; N=31.595338; O=0.7877119251503418; P=0.7877119251503418; Q=315.4033578243167; R=2838.630220418851; S=0.7328183953708449; T=0.732818395370845},
    @{A="ECs"; D="FAPs"; E=3; F=1; G=29.947775; H=89.84332500000001; I=0.9303126840830549; J=0.930312684083055; K=3; L=1; M=2.405620333333333; N=7.216861; O=0.1799255153355986; P=0.1799255153355986; Q=72.04297647809167; R=648.386788302825; S=0.1673869891068875; T=0.1673869891068876},
    @{A="ECs"; D="sCs"; E=3; F=1; G=29.947775; H=89.84332500000001; I=0.9303126840830549; J=0.930312684083055; K=2; L=0.6666666666666666; M=0.4326903333333333; N=1.298071; O=0.03236255951405961; P=0.03236255951405961; Q=12.95811274734167; R=116.623014726075; S=0.0301072996053224; T=0.03010729960532241},
    @{A="FAPs"; D="ECs"; E=3; F=1; G=0.9818753333333333; H=2.945626; I=0.03050146719708818; J=0.03050146719708818; K=2; L=0.6666666666666666; M=10.53177933333333; N=31.595338; O=0.7877119251503418; P=0.7877119251503418; Q=10.34089434350978; R=93.06804909158799; S=0.02402636944572833; T=0.02402636944572833},
    @{A="FAPs"; D="FAPs"; E=3; F=1; G=0.9818753333333333; H=2.945626; I=0.03050146719708818; J=0.03050146719708818; K=3; L=1; M=2.405620333333333; N=7.216861; O=0.1799255153355986; P=0.1799255153355986; Q=2.362019266665111; R=21.258173399986; S=0.005487992203927945; T=0.005487992203927946},
    @{A="FAPs"; D="sCs"; E=3; F=1; G=0.9818753333333333; H=2.945626; I=0.03050146719708818; J=0.03050146719708818; K=2; L=0.6666666666666666; M=0.4326903333333333; N=1.298071; O=0.03236255951405961; P=0.03236255951405961; Q=0.4248479652717778; R=3.823631687446; S=0.0009871055474319033; T=0.0009871055474319033},
    @{A="sCs"; D="ECs"; E=3; F=1; G=1.261435; H=3.784305; I=0.03918584871985675; J=0.03918584871985676; K=2; L=0.6666666666666666; M=10.53177933333333; N=31.595338; O=0.7877119251503418; P=0.7877119251503418; Q=13.28515506334333; R=119.56639557009; S=0.03086716033376842; T=0.03086716033376843},
    @{A="sCs"; D="FAPs"; E=3; F=1; G=1.261435; H=3.784305; I=0.03918584871985675; J=0.03918584871985676; K=3; L=1; M=2.405620333333333; N=7.216861; O=0.1799255153355986; P=0.1799255153355986; Q=3.034533685178333; R=27.310803166605; S=0.007050534024783031; T=0.007050534024783034},
    @{A="sCs"; D="sCs"; E=3; F=1; G=1.261435; H=3.784305; I=0.03918584871985675; J=0.03918584871985676; K=2; L=0.6666666666666666; M=0.4326903333333333; N=1.298071; O=0.03236255951405961; P=0.03236255951405961; Q=0.5458107306283333; R=4.912296575655; S=0.001268154361305301; T=0.001268154361305301}
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row["A"]
    $ws.Cells.Item($r, 2).Value = "Mmrn2"
    $ws.Cells.Item($r, 3).Value = "Clec14a"
    $ws.Cells.Item($r, 4).Value = $row["D"]
    $ws.Cells.Item($r, 5).Value = $row["E"]
    $ws.Cells.Item($r, 6).Value = $row["F"]
    $ws.Cells.Item($r, 7).Value = $row["G"]
    $ws.Cells.Item($r, 8).Value = $row["H"]
    $ws.Cells.Item($r, 9).Value = $row["I"]
    $ws.Cells.Item($r, 10).Value = $row["J"]
    $ws.Cells.Item($r, 11).Value = $row["K"]
    $ws.Cells.Item($r, 12).Value = $row["L"]
    $ws.Cells.Item($r, 13).Value = $row["M"]
    $ws.Cells.Item($r, 14).Value = $row["N"]
    $ws.Cells.Item($r, 15).Value = $row["O"]
    $ws.Cells.Item($r, 16).Value = $row["P"]
    $ws.Cells.Item($r, 17).Value = $row["Q"]
    $ws.Cells.Item($r, 18).Value = $row["R"]
    $ws.Cells.Item($r, 19).Value = $row["S"]
    $ws.Cells.Item($r, 20).Value = $row["T"]
}
